$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column H ("medium_light_green"), shifting columns I:L left to H:K
$ws.Columns.Item(8).Delete()

# Select the (now empty) column H, as the last user action
$ws.Range("H1:H1048576").Select()
